# "Generate Report for Handoff" - adds two new handed-off files
# (2f0ce313-1264-4909-81c7-cba41663550c and a57a5beb-f5f0-4cc4-851b-5b3be11ea5ea)
# to the localization status report, on all three sheets (Overview, zh-cn, de-de).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: Overview
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item(1)

# Clear any existing hyperlinks - we rebuild them all below since the
# runtime does not auto-shift hyperlink anchors when rows move.
$wsOverview.Hyperlinks.Delete()

# Row 2 (8b8c3d06) is unchanged. Rewrite rows 3-5.
$wsOverview.Range("A3").Value = "2f0ce313-1264-4909-81c7-cba41663550c.md"
$wsOverview.Range("B3").Value = "e2e\2f0ce313-1264-4909-81c7-cba41663550c.md"
$wsOverview.Range("C3").Value = ".md"
$wsOverview.Range("D3").Value = ""
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-08-19 22:47:17"

$wsOverview.Range("A4").Value = "a57a5beb-f5f0-4cc4-851b-5b3be11ea5ea.md"
$wsOverview.Range("B4").Value = "e2e\a57a5beb-f5f0-4cc4-851b-5b3be11ea5ea.md"
$wsOverview.Range("C4").Value = ".md"
$wsOverview.Range("D4").Value = ""
$wsOverview.Range("E4").Value = "Ready for handoff"
$wsOverview.Range("F4").Value = "Ready for handoff"
$wsOverview.Range("G4").Value = "2016-08-19 22:47:17"

$wsOverview.Range("A5").Value = "a6a8eb4c-617d-48b1-8154-a82f6da66c87.md"
$wsOverview.Range("B5").Value = "e2e\a6a8eb4c-617d-48b1-8154-a82f6da66c87.md"
$wsOverview.Range("C5").Value = ".md"
$wsOverview.Range("D5").Value = ""
$wsOverview.Range("E5").Value = "Ready for handoff"
$wsOverview.Range("F5").Value = "Ready for handoff"
$wsOverview.Range("G5").Value = "2016-08-19 22:46:11"

# Rebuild hyperlinks in order (B2..B5)
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4440f559ddc9459f83486cacb386569fb0d5201d/e2e/8b8c3d06-ff71-4797-83a0-86447618644c.md", "", "", "e2e\8b8c3d06-ff71-4797-83a0-86447618644c.md") | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2f0ce3131264490981c7cba41663550cabcdef0/e2e/2f0ce313-1264-4909-81c7-cba41663550c.md", "", "", "e2e\2f0ce313-1264-4909-81c7-cba41663550c.md") | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("B4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a57a5bebf5f04cc4851b5b3be11ea5eaabcdef0/e2e/a57a5beb-f5f0-4cc4-851b-5b3be11ea5ea.md", "", "", "e2e\a57a5beb-f5f0-4cc4-851b-5b3be11ea5ea.md") | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("B5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3ee1343e3bacd38b082c0640db725298a6e293a8/e2e/a6a8eb4c-617d-48b1-8154-a82f6da66c87.md", "", "", "e2e\a6a8eb4c-617d-48b1-8154-a82f6da66c87.md") | Out-Null

# Resize the Overview table (table3) and dimension to A1:G5
$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.Resize($wsOverview.Range("A1:G5"))

# ---------------------------------------------------------------------------
# Sheet 2: zh-cn
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item(2)
$wsZhCn.Hyperlinks.Delete()

# Row 2 (8b8c3d06) is unchanged. Rewrite rows 3-5.
$wsZhCn.Range("A3").Value = "2f0ce313-1264-4909-81c7-cba41663550c.md"
$wsZhCn.Range("B3").Value = ".md"
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("D3").Value = "e2e"
$wsZhCn.Range("E3").Value = "ht"
$wsZhCn.Range("F3").Value = "False"
$wsZhCn.Range("G3").Value = "2f0ce313-1264-4909-81c7-cba41663550c.6a6814d3fdbff2e0fa2289f007dc214b838fe235.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-08-19 22:47:13"
$wsZhCn.Range("I3").Value = ""
$wsZhCn.Range("J3").Value = ""
$wsZhCn.Range("K3").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("L3").Value = ""
$wsZhCn.Range("M3").Value = "True"
$wsZhCn.Range("N3").Value = ""
$wsZhCn.Range("O3").Value = "False"
$wsZhCn.Range("P3").Value = ""

$wsZhCn.Range("A4").Value = "a57a5beb-f5f0-4cc4-851b-5b3be11ea5ea.md"
$wsZhCn.Range("B4").Value = ".md"
$wsZhCn.Range("C4").Value = "Ready for handoff"
$wsZhCn.Range("D4").Value = "e2e"
$wsZhCn.Range("E4").Value = "ht"
$wsZhCn.Range("F4").Value = "False"
$wsZhCn.Range("G4").Value = "a57a5beb-f5f0-4cc4-851b-5b3be11ea5ea.e084d713dba2d996cfce88e53f12b0228de706a2.zh-cn.xlf"
$wsZhCn.Range("H4").Value = "2016-08-19 22:47:13"
$wsZhCn.Range("I4").Value = ""
$wsZhCn.Range("J4").Value = ""
$wsZhCn.Range("K4").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("L4").Value = ""
$wsZhCn.Range("M4").Value = "True"
$wsZhCn.Range("N4").Value = ""
$wsZhCn.Range("O4").Value = "False"
$wsZhCn.Range("P4").Value = ""

$wsZhCn.Range("A5").Value = "a6a8eb4c-617d-48b1-8154-a82f6da66c87.md"
$wsZhCn.Range("B5").Value = ".md"
$wsZhCn.Range("C5").Value = "Ready for handoff"
$wsZhCn.Range("D5").Value = "e2e"
$wsZhCn.Range("E5").Value = "ht"
$wsZhCn.Range("F5").Value = "False"
$wsZhCn.Range("G5").Value = "a6a8eb4c-617d-48b1-8154-a82f6da66c87.0c23cc8aa105ddcbe8118ae7ff150d474056f015.zh-cn.xlf"
$wsZhCn.Range("H5").Value = "2016-08-19 22:46:07"
$wsZhCn.Range("I5").Value = ""
$wsZhCn.Range("J5").Value = ""
$wsZhCn.Range("K5").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("L5").Value = ""
$wsZhCn.Range("M5").Value = "True"
$wsZhCn.Range("N5").Value = ""
$wsZhCn.Range("O5").Value = "False"
$wsZhCn.Range("P5").Value = ""

# Rebuild hyperlinks in order (A2, I2, A3, A4, A5)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4440f559ddc9459f83486cacb386569fb0d5201d/e2e/8b8c3d06-ff71-4797-83a0-86447618644c.md", "", "", "8b8c3d06-ff71-4797-83a0-86447618644c.md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/066aaad1ef678cb3962471db65444f640f1284e5/e2e/8b8c3d06-ff71-4797-83a0-86447618644c.md", "", "", "8b8c3d06-ff71-4797-83a0-86447618644c.md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2f0ce3131264490981c7cba41663550cabcdef0/e2e/2f0ce313-1264-4909-81c7-cba41663550c.md", "", "", "2f0ce313-1264-4909-81c7-cba41663550c.md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a57a5bebf5f04cc4851b5b3be11ea5eaabcdef0/e2e/a57a5beb-f5f0-4cc4-851b-5b3be11ea5ea.md", "", "", "a57a5beb-f5f0-4cc4-851b-5b3be11ea5ea.md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3ee1343e3bacd38b082c0640db725298a6e293a8/e2e/a6a8eb4c-617d-48b1-8154-a82f6da66c87.md", "", "", "a6a8eb4c-617d-48b1-8154-a82f6da66c87.md") | Out-Null

$loZhCn = $wsZhCn.ListObjects.Item(1)
$loZhCn.Resize($wsZhCn.Range("A1:P5"))

# ---------------------------------------------------------------------------
# Sheet 3: de-de
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item(3)
$wsDeDe.Hyperlinks.Delete()

# Row 2 (8b8c3d06) is unchanged. Rewrite rows 3-5.
$wsDeDe.Range("A3").Value = "2f0ce313-1264-4909-81c7-cba41663550c.md"
$wsDeDe.Range("B3").Value = ".md"
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("D3").Value = "e2e"
$wsDeDe.Range("E3").Value = "ht"
$wsDeDe.Range("F3").Value = "False"
$wsDeDe.Range("G3").Value = "2f0ce313-1264-4909-81c7-cba41663550c.6a6814d3fdbff2e0fa2289f007dc214b838fe235.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-08-19 22:47:17"
$wsDeDe.Range("I3").Value = ""
$wsDeDe.Range("J3").Value = ""
$wsDeDe.Range("K3").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("L3").Value = ""
$wsDeDe.Range("M3").Value = "True"
$wsDeDe.Range("N3").Value = ""
$wsDeDe.Range("O3").Value = "False"
$wsDeDe.Range("P3").Value = ""

$wsDeDe.Range("A4").Value = "a57a5beb-f5f0-4cc4-851b-5b3be11ea5ea.md"
$wsDeDe.Range("B4").Value = ".md"
$wsDeDe.Range("C4").Value = "Ready for handoff"
$wsDeDe.Range("D4").Value = "e2e"
$wsDeDe.Range("E4").Value = "ht"
$wsDeDe.Range("F4").Value = "False"
$wsDeDe.Range("G4").Value = "a57a5beb-f5f0-4cc4-851b-5b3be11ea5ea.e084d713dba2d996cfce88e53f12b0228de706a2.de-de.xlf"
$wsDeDe.Range("H4").Value = "2016-08-19 22:47:17"
$wsDeDe.Range("I4").Value = ""
$wsDeDe.Range("J4").Value = ""
$wsDeDe.Range("K4").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("L4").Value = ""
$wsDeDe.Range("M4").Value = "True"
$wsDeDe.Range("N4").Value = ""
$wsDeDe.Range("O4").Value = "False"
$wsDeDe.Range("P4").Value = ""

$wsDeDe.Range("A5").Value = "a6a8eb4c-617d-48b1-8154-a82f6da66c87.md"
$wsDeDe.Range("B5").Value = ".md"
$wsDeDe.Range("C5").Value = "Ready for handoff"
$wsDeDe.Range("D5").Value = "e2e"
$wsDeDe.Range("E5").Value = "ht"
$wsDeDe.Range("F5").Value = "False"
$wsDeDe.Range("G5").Value = "a6a8eb4c-617d-48b1-8154-a82f6da66c87.0c23cc8aa105ddcbe8118ae7ff150d474056f015.de-de.xlf"
$wsDeDe.Range("H5").Value = "2016-08-19 22:46:11"
$wsDeDe.Range("I5").Value = ""
$wsDeDe.Range("J5").Value = ""
$wsDeDe.Range("K5").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("L5").Value = ""
$wsDeDe.Range("M5").Value = "True"
$wsDeDe.Range("N5").Value = ""
$wsDeDe.Range("O5").Value = "False"
$wsDeDe.Range("P5").Value = ""

# Rebuild hyperlinks in order (A2, I2, A3, A4, A5)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4440f559ddc9459f83486cacb386569fb0d5201d/e2e/8b8c3d06-ff71-4797-83a0-86447618644c.md", "", "", "8b8c3d06-ff71-4797-83a0-86447618644c.md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/bf94d71aa848ef0dbca62a0cd9b8bff482f32e93/e2e/8b8c3d06-ff71-4797-83a0-86447618644c.md", "", "", "8b8c3d06-ff71-4797-83a0-86447618644c.md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2f0ce3131264490981c7cba41663550cabcdef0/e2e/2f0ce313-1264-4909-81c7-cba41663550c.md", "", "", "2f0ce313-1264-4909-81c7-cba41663550c.md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a57a5bebf5f04cc4851b5b3be11ea5eaabcdef0/e2e/a57a5beb-f5f0-4cc4-851b-5b3be11ea5ea.md", "", "", "a57a5beb-f5f0-4cc4-851b-5b3be11ea5ea.md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3ee1343e3bacd38b082c0640db725298a6e293a8/e2e/a6a8eb4c-617d-48b1-8154-a82f6da66c87.md", "", "", "a6a8eb4c-617d-48b1-8154-a82f6da66c87.md") | Out-Null

$loDeDe = $wsDeDe.ListObjects.Item(1)
$loDeDe.Resize($wsDeDe.Range("A1:P5"))
